$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109:217 down to 110:218
$ws.Rows(109).Insert()

# Populate the newly inserted row 109 with its data
$ws.Range("A109").Value = 5
$ws.Range("B109").Value = "Macroferia Regional de Talca"
$ws.Range("C109").Value = "Maule"
$ws.Range("D109").Value = 44566
$ws.Range("D109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E109").Value = 7
$ws.Range("F109").Value = 100112006
$ws.Range("G109").Value = "Repollo"
$ws.Range("H109").Value = "Crespo record"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 700
$ws.Range("L109").Value = 700
$ws.Range("M109").Value = 700
$ws.Range("N109").Value = "$/unidad"
$ws.Range("O109").Value = "Región del Maule"
$ws.Range("P109").Value = 700
$ws.Range("Q109").Value = 1
$ws.Range("R109").Value = "Hortaliza"
